# Apply the "Trade #11 closed" update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet: update aggregate stats (now 11 trades total, 8 losing)
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B4").Value = -0.34   # Total P&L $
$summary.Range("B5").Value = -0.62   # Total P&L %
$summary.Range("B6").Value = 11      # Total Trades
$summary.Range("B8").Value = 8       # Losing Trades
$summary.Range("B9").Value = 27.27   # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet: update the MarketMaking strategy row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 11       # Trades
$status.Range("E5").Value = -0.34    # P&L $
$status.Range("G5").Value = 27.27    # Win Rate %

# ---------------------------------------------------------------------------
# 3) All Trades + MarketMaking sheets: append the newly closed trade (row 12)
# ---------------------------------------------------------------------------
$newTradeRow = @{
    A = 11
    B = "2026-02-17"
    C = "20:03:02"
    D = "MarketMaking"
    E = "UP"
    F = 0.76
    G = 0.750988
    H = "CLOSED"
    I = -1.1858
    J = -0.01
    K = 99.67
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A12").Value = $newTradeRow.A
    # Force text so the "2026-02-17" date-like string isn't auto-converted
    # into a date serial number (matches the existing rows' storage as text).
    $ws.Range("B12").NumberFormat = "@"
    $ws.Range("B12").Value = $newTradeRow.B
    $ws.Range("C12").Value = $newTradeRow.C
    $ws.Range("D12").Value = $newTradeRow.D
    $ws.Range("E12").Value = $newTradeRow.E
    $ws.Range("F12").Value = $newTradeRow.F
    $ws.Range("G12").Value = $newTradeRow.G
    $ws.Range("H12").Value = $newTradeRow.H
    $ws.Range("I12").Value = $newTradeRow.I
    $ws.Range("J12").Value = $newTradeRow.J
    $ws.Range("K12").Value = $newTradeRow.K
    $ws.Range("L12").Value = $newTradeRow.L
    $ws.Range("M12").Value = $newTradeRow.M
    $ws.Range("N12").Value = $newTradeRow.N
    $ws.Range("O12").Value = $newTradeRow.O
    $ws.Range("P12").Value = $newTradeRow.P
    $ws.Range("Q12").Value = $newTradeRow.Q
}
